# The presentation had three tables (on slides 14, 15 and 16) that were
# using a plain, locally-defined table style. The author re-styled all
# three with PowerPoint's built-in "Themed Style 1 - Accent 1" table
# style ({5EFFE562-3217-48EB-B695-77BA6CED0536}), replacing the former
# custom style ({1FADA30B-C856-48E8-9ED7-D0541508CE6B}).
#
# Table styles can't be changed by assigning to Table.Style (that
# property is read-only in the object model) - PowerPoint exposes the
# mutator as Table.ApplyStyle("{GUID}").

$p = $ppt.ActivePresentation

$newStyleId = "{5EFFE562-3217-48EB-B695-77BA6CED0536}"
$slideIndexesWithRestyledTables = 14, 15, 16

foreach ($slideIndex in $slideIndexesWithRestyledTables) {
    $slide = $p.Slides.Item($slideIndex)

    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
